$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.053.87"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "  -0.01%  "
$ws.Range("D3").Value = "'1.835.51"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "  +0.37%  "
$ws.Range("D4").Value = "'0.9994"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "  +0.07%  "
$ws.Range("D5").Value = "'243.97"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +1.30%  "
$ws.Range("E6").Value = "  +2.20%  "
$ws.Range("E7").Value = "  +0.04%  "
$ws.Range("D8").Value = "'0.07580"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "  +3.15%  "
$ws.Range("D9").Value = "'0.2949"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +1.17%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.07740"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "  +0.72%  "
$ws.Range("D12").Value = "'1.840.67"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +0.68%  "
$ws.Range("D13").Value = "'4.995"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +0.57%  "
$ws.Range("D14").Value = "'0.6709"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "  +1.36%  "
$ws.Range("D15").Value = "'83.27"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +1.72%  "
$ws.Range("D16").Value = "'0.000009842"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "  +9.67%  "
$ws.Range("D17").Value = "'6.121"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.63%  "
$ws.Range("D18").Value = "'29.091.27"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +0.12%  "
$ws.Range("D19").Value = "'12.56"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +1.62%  "
$ws.Range("D20").Value = "'226.94"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "  +0.56%  "
$ws.Range("E21").Value = "  -0.05%  "
$ws.Range("D22").Value = "'7.236"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +1.42%  "
$ws.Range("E23").Value = "  +0.06%  "
$ws.Range("D24").Value = "'160.51"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +0.60%  "
$ws.Range("D25").Value = "'0.1406"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +3.83%  "
$ws.Range("D26").Value = "'8.550"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +1.42%  "
$ws.Range("D27").Value = "'17.97"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "  +1.02%  "
$ws.Range("D28").Value = "'1.503"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +0.53%  "
$ws.Range("D29").Value = "'4.123"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "  +1.81%  "
$ws.Range("D30").Value = "'4.056"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +0.37%  "
$ws.Range("D31").Value = "'1.201"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "'0.05362"
$ws.Range("D32").ClearFormats()
$ws.Range("E32").Value = "  +1.98%  "
$ws.Range("D33").Value = "'1.862"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").Value = "'0.7502"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +2.80%  "
$ws.Range("E35").Value = "  -0.71%  "
$ws.Range("D36").Value = "'2.671"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +0.92%  "
$ws.Range("D37").Value = "'1.250.89"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "  -3.39%  "
$ws.Range("D38").Value = "'0.01798"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "  +0.83%  "
$ws.Range("D39").Value = "'2.761"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +0.45%  "
$ws.Range("D40").Value = "'6.602"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +4.99%  "
$ws.Range("D41").Value = "'0.9055"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "  +0.58%  "
$ws.Range("E42").Value = "  +0.17%  "
$ws.Range("D43").Value = "'102.60"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +0.65%  "
$ws.Range("D44").Value = "'1.985.00"
$ws.Range("D44").ClearFormats()
$ws.Range("E44").Value = "  +0.57%  "
$ws.Range("D45").Value = "'0.00000000124"
$ws.Range("D45").ClearFormats()
$ws.Range("E45").Value = "  +3.59%  "
$ws.Range("E46").Value = "  +1.61%  "
$ws.Range("D47").Value = "'0.5116"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +0.03%  "
$ws.Range("D48").Value = "'0.4097"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("D49").Value = "'9.089"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "  +3.20%  "
$ws.Range("D50").Value = "'0.05802"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.25%  "
$ws.Range("D51").Value = "'6.770"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "  +1.65%  "
